# Auto-generated script applying the Sagittarius_Profits.xlsx data refresh
# (scheduled runner pulled updated Leve profit figures for ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3560
$ws.Range("I2").Value = 620
$ws.Range("J2").Value = 7235
$ws.Range("K2").Value = 620
$ws.Range("L2").Value = 7235
$ws.Range("M2").Value = -507
$ws.Range("N2").Value = -7461

$ws.Range("H116").Value = 8502.866
$ws.Range("I116").Value = 8510.200000000001
$ws.Range("K116").Value = 8510.200000000001
$ws.Range("M116").Value = -5068.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1184
$ws.Range("I2").Value = 670.8
$ws.Range("K2").Value = 670.8
$ws.Range("M2").Value = -557.8

$ws.Range("H10").Value = 305
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 305
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 305
$ws.Range("N10").Value = -645
$ws.Range("M10").ClearContents()

$ws.Range("H74").Value = 2764.7896
$ws.Range("I74").Value = 2058.077
$ws.Range("J74").Value = 4296
$ws.Range("K74").Value = 2058.077
$ws.Range("L74").Value = 4296
$ws.Range("M74").Value = -1184.077
$ws.Range("N74").Value = -6044

$ws.Range("H77").Value = 2764.7896
$ws.Range("I77").Value = 2058.077
$ws.Range("J77").Value = 4296
$ws.Range("K77").Value = 10290.385
$ws.Range("L77").Value = 21480
$ws.Range("M77").Value = -5922.385000000002
$ws.Range("N77").Value = -30216

$ws.Range("H116").Value = 1184
$ws.Range("I116").Value = 670.8
$ws.Range("K116").Value = 670.8
$ws.Range("M116").Value = 1623.2

$ws.Range("H122").Value = 2274.3
$ws.Range("I122").Value = 1773.5
$ws.Range("J122").Value = 2608.1667
$ws.Range("K122").Value = 5320.5
$ws.Range("L122").Value = 7824.500100000001
$ws.Range("M122").Value = -2870.5
$ws.Range("N122").Value = -12724.5001

$ws.Range("H132").Value = 2070
$ws.Range("I132").Value = 1962.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5887.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3357.5
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1184
$ws.Range("I3").Value = 670.8
$ws.Range("K3").Value = 670.8
$ws.Range("M3").Value = -556.8

$ws.Range("H20").Value = 1096.75
$ws.Range("I20").Value = 1162.6666
$ws.Range("J20").Value = 899
$ws.Range("K20").Value = 1162.6666
$ws.Range("L20").Value = 899
$ws.Range("M20").Value = -915.6666
$ws.Range("N20").Value = -1393

$ws.Range("H86").Value = 1392.5
$ws.Range("I86").Value = 1392.5
$ws.Range("K86").Value = 1392.5
$ws.Range("M86").Value = -269.5

$ws.Range("H89").Value = 1392.5
$ws.Range("I89").Value = 1392.5
$ws.Range("K89").Value = 6962.5
$ws.Range("M89").Value = -1346.5

$ws.Range("H134").Value = 2137.9
$ws.Range("I134").Value = 2264.4443
$ws.Range("K134").Value = 6793.3329
$ws.Range("M134").Value = -4258.3329

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 3468
$ws.Range("I5").Value = 3468
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3468
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3356
$ws.Range("N5").ClearContents()

$ws.Range("H16").Value = 1645.3158
$ws.Range("I16").Value = 1276.6
$ws.Range("K16").Value = 1276.6
$ws.Range("M16").Value = -989.5999999999999

$ws.Range("H58").Value = 2259.3076
$ws.Range("I58").Value = 2540.8333
$ws.Range("K58").Value = 2540.8333
$ws.Range("M58").Value = -2337.8333

$ws.Range("H105").Value = 3545.4075
$ws.Range("I105").Value = 2748.8572
$ws.Range("J105").Value = 4403.231
$ws.Range("K105").Value = 2748.8572
$ws.Range("L105").Value = 4403.231
$ws.Range("M105").Value = -1001.8572
$ws.Range("N105").Value = -7897.231

$ws.Range("H113").Value = 1645.3158
$ws.Range("I113").Value = 1276.6
$ws.Range("K113").Value = 1276.6
$ws.Range("M113").Value = 893.4000000000001

$ws.Range("H136").Value = 2259.3076
$ws.Range("I136").Value = 2540.8333
$ws.Range("K136").Value = 7622.499899999999
$ws.Range("M136").Value = -5072.499899999999

$ws.Range("H141").Value = 411396.9
$ws.Range("J141").Value = 411396.9
$ws.Range("L141").Value = 411396.9
$ws.Range("N141").Value = -421756.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 22
$ws.Range("J16").Value = 22
$ws.Range("L16").Value = 66
$ws.Range("N16").Value = -412

$ws.Range("H34").Value = 240
$ws.Range("I34").Value = 240
$ws.Range("K34").Value = 720
$ws.Range("M34").Value = -636

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H107").Value = 1019.4
$ws.Range("J107").Value = 1997.5
$ws.Range("L107").Value = 5992.5
$ws.Range("N107").Value = -9832.5

$ws.Range("H123").Value = 10734.143
$ws.Range("I123").Value = 4996.6665
$ws.Range("J123").Value = 15037.25
$ws.Range("K123").Value = 14989.9995
$ws.Range("L123").Value = 45111.75
$ws.Range("M123").Value = -12539.9995
$ws.Range("N123").Value = -50011.75

$ws.Range("H134").Value = 13402.7
$ws.Range("I134").Value = 1798.3334
$ws.Range("K134").Value = 5395.0002
$ws.Range("M134").Value = -325.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H43").Value = 9011.333000000001
$ws.Range("I43").Value = 9011.333000000001
$ws.Range("K43").Value = 9011.333000000001
$ws.Range("M43").Value = -8860.333000000001

$ws.Range("H70").Value = 6750
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6750
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6750
$ws.Range("N70").Value = -7290
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 6750
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6750
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6750
$ws.Range("N73").Value = -8622
$ws.Range("M73").ClearContents()

$ws.Range("H113").Value = 1724.75
$ws.Range("I113").Value = 1633
$ws.Range("K113").Value = 1633
$ws.Range("M113").Value = 537

$ws.Range("H122").Value = 4659
$ws.Range("I122").Value = 3035.6667
$ws.Range("K122").Value = 9107.000100000001
$ws.Range("M122").Value = -6657.000100000001

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2912
$ws.Range("I132").Value = 2912
$ws.Range("K132").Value = 8736
$ws.Range("M132").Value = -6206

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H25").Value = 1600
$ws.Range("I25").Value = 1600
$ws.Range("K25").Value = 1600
$ws.Range("M25").Value = -1370

$ws.Range("H61").Value = 3041.762
$ws.Range("I61").Value = 2443.95
$ws.Range("K61").Value = 2443.95
$ws.Range("M61").Value = -2241.95

$ws.Range("H113").Value = 3041.762
$ws.Range("I113").Value = 2443.95
$ws.Range("K113").Value = 2443.95
$ws.Range("M113").Value = -273.9499999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H23").Value = 990
$ws.Range("I23").Value = 990
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 990
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -761
$ws.Range("N23").ClearContents()

$ws.Range("H75").Value = 89000
$ws.Range("I75").Value = 89000
$ws.Range("K75").Value = 89000
$ws.Range("M75").Value = -88064

$ws.Range("H78").Value = 89000
$ws.Range("I78").Value = 89000
$ws.Range("K78").Value = 267000
$ws.Range("M78").Value = -262320

$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

$ws.Range("H122").Value = 644
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 5580.4165
$ws.Range("I126").Value = 3079.8333
$ws.Range("J126").Value = 8081
$ws.Range("K126").Value = 9239.499899999999
$ws.Range("L126").Value = 24243
$ws.Range("M126").Value = -6769.499899999999
$ws.Range("N126").Value = -29183

$ws.Range("H132").Value = 1799
$ws.Range("I132").Value = 1799
$ws.Range("K132").Value = 5397
$ws.Range("M132").Value = -2867

$ws.Range("H136").Value = 2430.25
$ws.Range("I136").Value = 2430.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7290.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4740.75
$ws.Range("N136").ClearContents()
